$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now becomes the values previously on row 6 (with M unchanged)
$ws.Range("D2").Value = 44263
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("S2").Value = 1194

# Row 3: now becomes the values previously on row 2
$ws.Range("D3").Value = 44307
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("S3").Value = 1083

# Row 6: now becomes the values previously on row 3
$ws.Range("D6").Value = 44323
$ws.Range("M6").Value = 270
$ws.Range("Q6").Value = "$/bandeja 18 kilos"
